$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp and country name swaps (column A) ---
$ws.Range("A1").Value = 'Datos actualizados a 17 de Octubre de 2020 a las 14:55'
$ws.Range("A60").Value = 'Austria'
$ws.Range("A61").Value = 'Armenia'
$ws.Range("A62").Value = 'Uzbekistan'
$ws.Range("A78").Value = 'Dinamarca'
$ws.Range("A79").Value = 'Tunez'
$ws.Range("A80").Value = 'Jordania'
$ws.Range("A81").Value = 'Bosnia y Herzegovina'
$ws.Range("A82").Value = 'Birmania'
$ws.Range("A216").Value = 'Islas Malvinas'
$ws.Range("A217").Value = 'Montserrat'

# --- Update numeric data cells (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 8291779
$ws.Range("C4").Value = 3501
$ws.Range("D4").Value = 5396689
$ws.Range("E4").Value = 2671395
$ws.Range("G4").Value = 51
$ws.Range("H4").Value = 223695
$ws.Range("B5").Value = 7434630
$ws.Range("C5").Value = 3995
$ws.Range("E5").Value = 796973
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 113062
$ws.Range("B18").Value = 423524
$ws.Range("C18").Value = 3221
$ws.Range("D18").Value = 357291
$ws.Range("E18").Value = 56035
$ws.Range("G18").Value = 56
$ws.Range("H18").Value = 10198
$ws.Range("D22").Value = 290000
$ws.Range("E22").Value = 56956
$ws.Range("B25").Value = 341854
$ws.Range("C25").Value = 359
$ws.Range("D25").Value = 328165
$ws.Range("E25").Value = 8524
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 5165
$ws.Range("B29").Value = 220052
$ws.Range("C29").Value = 8114
$ws.Range("G29").Value = 29
$ws.Range("H29").Value = 6737
$ws.Range("B42").Value = 115483
$ws.Range("C42").Value = 739
$ws.Range("D42").Value = 107108
$ws.Range("E42").Value = 7681
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 694
$ws.Range("B60").Value = 63134
$ws.Range("C60").Value = 1747
$ws.Range("D60").Value = 48771
$ws.Range("E60").Value = 13474
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 889
$ws.Range("B61").Value = 63000
$ws.Range("C61").Value = 1540
$ws.Range("D61").Value = 47925
$ws.Range("E61").Value = 14008
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 1067
$ws.Range("B62").Value = 62809
$ws.Range("C62").Value = 221
$ws.Range("D62").Value = 59756
$ws.Range("E62").Value = 2531
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 522
$ws.Range("B72").Value = 46746
$ws.Range("C72").Value = 312
$ws.Range("D72").Value = 40162
$ws.Range("E72").Value = 6182
$ws.Range("B74").Value = 44317
$ws.Range("C74").Value = 528
$ws.Range("D74").Value = 39903
$ws.Range("E74").Value = 3791
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 623
$ws.Range("B78").Value = 34941
$ws.Range("C78").Value = 500
$ws.Range("D78").Value = 28917
$ws.Range("E78").Value = 5345
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 679
$ws.Range("B79").Value = 34790
$ws.Range("D79").Value = 5032
$ws.Range("E79").Value = 29246
$ws.Range("H79").Value = 512
$ws.Range("B80").Value = 34548
$ws.Range("D80").Value = 6692
$ws.Range("E80").Value = 27546
$ws.Range("H80").Value = 310
$ws.Range("B81").Value = 33561
$ws.Range("C81").Value = 716
$ws.Range("D81").Value = 24773
$ws.Range("E81").Value = 7807
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 981
$ws.Range("B82").Value = 33488
$ws.Range("D82").Value = 15477
$ws.Range("E82").Value = 17212
$ws.Range("H82").Value = 799
$ws.Range("B99").Value = 15392
$ws.Range("C99").Value = 24
$ws.Range("D99").Value = 13756
$ws.Range("E99").Value = 1319
$ws.Range("B111").Value = 10455
$ws.Range("C111").Value = 41
$ws.Range("D111").Value = 9457
$ws.Range("E111").Value = 918
$ws.Range("B143").Value = 3998
$ws.Range("C143").Value = 69
$ws.Range("D143").Value = 2745
$ws.Range("E143").Value = 1242
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
